$d = $word.ActiveDocument

function Insert-LineBreak($findText, $replaceText) {
    $result = $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $result) {
        throw "Find.Execute did not find text: $findText"
    }
}

# Paragraph: Objetivos
Insert-LineBreak "Propiciar ao aluno uma visão básica sobre os principais métodos de determinação teórica da estrutura eletrônica dos materiais, com enfoque em sólidos cristalinos, mas também em materiais bidimensionais e nanoestruturados.O pr" "Propiciar ao aluno uma visão básica sobre os principais métodos de determinação teórica da estrutura eletrônica dos materiais, com enfoque em sólidos cristalinos, mas também em materiais bidimensionais e nanoestruturados.^lO pr"
Insert-LineBreak "O principal método de cálculo a ser empregado no curso será a Teoria do Funcional da Densidade(Den" "O principal método de cálculo a ser empregado no curso será a Teoria do Funcional da Densidade^l(Den"

# Paragraph: Programa (detailed)
Insert-LineBreak "Revisão de mecânica quânticao Eq" "Revisão de mecânica quântica^lo Eq"
Insert-LineBreak "o Equação de Schrödingero Át" "o Equação de Schrödinger^lo Át"
Insert-LineBreak "o Átomo do hidrogênio e orbitais atômicoso No" "o Átomo do hidrogênio e orbitais atômicos^lo No"
Insert-LineBreak "o Notação de Diraco Pr" "o Notação de Dirac^lo Pr"
Insert-LineBreak "o Princípio variacionalo Co" "o Princípio variacional^lo Co"
Insert-LineBreak "o Combinação linear de orbitais atômicosRevi" "o Combinação linear de orbitais atômicos^lRevi"
Insert-LineBreak "Revisão de física do estado sólidoo Es" "Revisão de física do estado sólido^lo Es"
Insert-LineBreak "o Espaço direto e recíprocoo Te" "o Espaço direto e recíproco^lo Te"
Insert-LineBreak "o Teorema de Blocho Zo" "o Teorema de Bloch^lo Zo"
Insert-LineBreak "o Zona de Brillouino Ba" "o Zona de Brillouin^lo Ba"
Insert-LineBreak "o Bandas de energia e densidade de estadoso En" "o Bandas de energia e densidade de estados^lo En"
Insert-LineBreak "o Energia de Fermi e superfície de Fermio Ap" "o Energia de Fermi e superfície de Fermi^lo Ap"
Insert-LineBreak "o Aproximação de elétrons livresMéto" "o Aproximação de elétrons livres^lMéto"
Insert-LineBreak "Método de Hartree-Focko De" "Método de Hartree-Fock^lo De"
Insert-LineBreak "o Determinantes de Slatero Eq" "o Determinantes de Slater^lo Eq"
Insert-LineBreak "o Equação de Hartree-Focko Po" "o Equação de Hartree-Fock^lo Po"
Insert-LineBreak "o Potencial de troca e correlaçãoo Al" "o Potencial de troca e correlação^lo Al"
Insert-LineBreak "o Algoritmo autoconsistenteTeor" "o Algoritmo autoconsistente^lTeor"
Insert-LineBreak "Teoria do funcional da densidadeo Te" "Teoria do funcional da densidade^lo Te"
Insert-LineBreak "o Teoremas de Hohenberg-Kohno Eq" "o Teoremas de Hohenberg-Kohn^lo Eq"
Insert-LineBreak "o Equações de Kohn-Shamo Fu" "o Equações de Kohn-Sham^lo Fu"
Insert-LineBreak "o Funcionais de troca e correlação: LDA, GGA, etc.Méto" "o Funcionais de troca e correlação: LDA, GGA, etc.^lMéto"
Insert-LineBreak "Métodos de ondas planas e pseudo-potenciaiso Ba" "Métodos de ondas planas e pseudo-potenciais^lo Ba"
Insert-LineBreak "o Bases de ondas planaso Ps" "o Bases de ondas planas^lo Ps"
Insert-LineBreak "o Pseudo-potenciaiso Ba" "o Pseudo-potenciais^lo Ba"
Insert-LineBreak "o Bases de ondas planas aumentadas e linearizadaso Mé" "o Bases de ondas planas aumentadas e linearizadas^lo Mé"
Insert-LineBreak "o Método FP-LAPWCódi" "o Método FP-LAPW^lCódi"
Insert-LineBreak "Códigos computacionaiso Qu" "Códigos computacionais^lo Qu"
Insert-LineBreak "o Quantum Espressoo El" "o Quantum Espresso^lo El"
Insert-LineBreak "o Elko Wi" "o Elk^lo Wi"
Insert-LineBreak "o Wien2ko VA" "o Wien2k^lo VA"

Write-Host "Done"